{"js": "// Revert \"more rule fixes, general rename and reorder\":\n// In the proof table, the row whose \"Line\" value is \"4\" should have its\n// \"First Segment\" (was \"2\") and \"Second Segment\" (was \"3\") cells cleared\n// back to empty, matching the other rows that don't cite two segments.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every cell's text so we can locate the header columns and the\n// target data row by content instead of relying on fixed indexes.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\n// Word cell text always ends with a trailing control character (tab\n// between cells / cell-mark on the last one) - strip it for comparisons.\nfunction cellText(cell) {\n  return cell.body.text.replace(/[\\s\\u0000-\\u001f]+$/, \"\");\n}\n\n// Row 0 is the header: find the \"First Segment\" / \"Second Segment\" columns.\nconst headerCells = rows.items[0].cells.items;\nlet firstSegCol = -1;\nlet secondSegCol = -1;\nfor (let c = 0; c < headerCells.length; c++) {\n  const label = cellText(headerCells[c]);\n  if (label === \"First Segment\") firstSegCol = c;\n  if (label === \"Second Segment\") secondSegCol = c;\n}\n\n// Find the data row whose first column (\"Line\") reads \"4\".\nlet targetRow = -1;\nfor (let r = 1; r < rows.items.length; r++) {\n  if (cellText(rows.items[r].cells.items[0]) === \"4\") {\n    targetRow = r;\n    break;\n  }\n}\n\nif (targetRow !== -1 && firstSegCol !== -1 && secondSegCol !== -1) {\n  const cellsToClear = [\n    rows.items[targetRow].cells.items[firstSegCol],\n    rows.items[targetRow].cells.items[secondSegCol],\n  ];\n\n  for (const cell of cellsToClear) {\n    const text = cellText(cell);\n    if (text.length === 0) {\n      continue; // already empty\n    }\n    // Search for just the visible text (excludes the trailing cell mark)\n    // so removing it leaves a plain empty paragraph, instead of wiping\n    // out the whole paragraph/cell content.\n    const bodyRange = cell.body.getRange();\n    const found = bodyRange.search(text);\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length > 0) {\n      found.items[0].insertText(\"\", Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Revert \"more rule fixes, general rename and reorder\":\n# In the proof table, the row whose \"Line\" value is \"4\" should have its\n# \"First Segment\" (was \"2\") and \"Second Segment\" (was \"3\") cells cleared\n# back to empty, matching the other rows that don't cite two segments.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n# A Word cell's Range.Text always carries a trailing control character\n# (tab between cells / bell on the last cell of a row) - strip it so we\n# can compare against the plain visible text.\nfunction Get-CellText($cell) {\n    return ($cell.Range.Text -replace \"[\\s\\x00-\\x1f]+$\", \"\")\n}\n\n# Locate the \"First Segment\" / \"Second Segment\" columns from the header row.\n$firstSegCol = -1\n$secondSegCol = -1\nfor ($c = 1; $c -le $colCount; $c++) {\n    $header = Get-CellText $table.Cell(1, $c)\n    if ($header -eq \"First Segment\") { $firstSegCol = $c }\n    if ($header -eq \"Second Segment\") { $secondSegCol = $c }\n}\n\n# Locate the data row whose \"Line\" column reads \"4\".\n$targetRow = -1\nfor ($r = 2; $r -le $rowCount; $r++) {\n    $lineText = Get-CellText $table.Cell($r, 1)\n    if ($lineText -eq \"4\") {\n        $targetRow = $r\n        break\n    }\n}\n\nif ($targetRow -ne -1 -and $firstSegCol -ne -1 -and $secondSegCol -ne -1) {\n    $cols = @($firstSegCol, $secondSegCol)\n    foreach ($col in $cols) {\n        $cell = $table.Cell($targetRow, $col)\n        $cellRange = $cell.Range\n        # Exclude the trailing end-of-cell marker so only the visible\n        # text is removed and the (now empty) paragraph is preserved.\n        $clearRange = $d.Range($cellRange.Start, $cellRange.End - 1)\n        $clearRange.Delete()\n    }\n}\n"}
